$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 10.1579178889515
$ws.Range("C2").Value = 19.598188785896209
$ws.Range("D2").Value = 54.960389128487662
$ws.Range("E2").Value = 27.878270548128647
$ws.Range("F2").Value = 33.829393954056201
$ws.Range("G2").Value = 41.324486592780211
$ws.Range("H2").Value = 15.746230512527513
$ws.Range("I2").Value = 24.38190647413434
$ws.Range("J2").Value = 31.70431930037735
$ws.Range("K2").Value = 20.604382251225044
$ws.Range("L2").Value = 24.322952875451943
$ws.Range("M2").Value = 47.843257123483276
$ws.Range("N2").Value = 17.312042827706232
$ws.Range("O2").Value = 44.663834761267047
$ws.Range("P2").Value = 11.284944697953476
$ws.Range("Q2").Value = 19.098131072913013
$ws.Range("R2").Value = 44.342436475661259
$ws.Range("T2").Value = 31.602233497446964
$ws.Range("U2").Value = 44.310338944682002
$ws.Range("V2").Value = 50.738934451990168
$ws.Range("W2").Value = 34.868621084734421
$ws.Range("X2").Value = 29.051183153235279
$ws.Range("Y2").Value = 46.873127869532603
$ws.Range("Z2").Value = 74.102129055464218
$ws.Range("AA2").Value = 58.137537219299716
$ws.Range("AB2").Value = 82.117012877220645
$ws.Range("AC2").Value = 47.226731427741825
$ws.Range("AD2").Value = 68.118386208681187
$ws.Range("AE2").Value = 37.334934363586754
$ws.Range("AF2").Value = 50.683503696834187
$ws.Range("AG2").Value = 43.103091125518958
$ws.Range("AH2").Value = 51.496437210824041
$ws.Range("AI2").Value = 71.501268927193323
$ws.Range("AJ2").Value = 43.110371741536255
$ws.Range("AK2").Value = 31.446994049177963
$ws.Range("AL2").Value = 26.34763984738133
$ws.Range("AM2").Value = 51.934247990362337
$ws.Range("AN2").Value = 48.870566677367158
$ws.Range("AO2").Value = 51.514734871571356
$ws.Range("AP2").Value = 42.763013751717352
$ws.Range("AQ2").Value = 17.959769131532035
$ws.Range("AR2").Value = 40.147522111484527
$ws.Range("AS2").Value = 49.212317284231553
$ws.Range("AT2").Value = 68.981116870872142
$ws.Range("AU2").Value = 26.844812428814741
$ws.Range("AV2").Value = 69.087466343055254
$ws.Range("AW2").Value = 61.535474026371006
$ws.Range("AX2").Value = 60.975739966022303
$ws.Range("AY2").Value = 61.042352709855827
$ws.Range("B3").Value = 26.86979334674643
$ws.Range("C3").Value = 34.915126866031834
$ws.Range("D3").Value = 48.851432707193773
$ws.Range("E3").Value = 183.28514142988689
$ws.Range("F3").Value = 25.143259564549318
$ws.Range("G3").Value = 35.906172158586244
$ws.Range("H3").Value = 26.237936367263348
$ws.Range("I3").Value = 15.480859712581388
$ws.Range("J3").Value = 45.489901724455379
$ws.Range("K3").Value = 31.557407905252216
$ws.Range("L3").Value = 29.155936108453172
$ws.Range("M3").Value = 35.950924922645264
$ws.Range("N3").Value = 30.333193715840395
$ws.Range("O3").Value = 68.291713134781546
$ws.Range("P3").Value = 59.724146844318859
$ws.Range("Q3").Value = 17.112738135193769
$ws.Range("R3").Value = 40.919423938942096
$ws.Range("S3").Value = 39.23765354394515
$ws.Range("T3").Value = 35.286204842344496
$ws.Range("U3").Value = 39.885976679604532
$ws.Range("V3").Value = 28.294614371481259
$ws.Range("W3").Value = 39.449478634693143
$ws.Range("X3").Value = 31.10636744546003
$ws.Range("Y3").Value = 73.157759855144505
$ws.Range("Z3").Value = 73.691744684544176
$ws.Range("AA3").Value = 58.505369020645169
$ws.Range("AB3").Value = 57.654260474615846
$ws.Range("AC3").Value = 35.298459469591322
$ws.Range("AD3").Value = 56.708382321554907
$ws.Range("AE3").Value = 52.994730892147466
$ws.Range("AF3").Value = 61.831623111663717
$ws.Range("AG3").Value = 35.101390433806287
$ws.Range("AH3").Value = 34.139991730003572
$ws.Range("AI3").Value = 55.246575236677209
$ws.Range("AJ3").Value = 59.414680685436416
$ws.Range("AK3").Value = 40.872347734999551
$ws.Range("AL3").Value = 46.019329743752188
$ws.Range("AM3").Value = 69.12796027587531
$ws.Range("AN3").Value = 47.032783948673362
$ws.Range("AO3").Value = 87.695087615686973
$ws.Range("AP3").Value = 50.28015910235689
$ws.Range("AQ3").Value = 28.062484626069789
$ws.Range("AR3").Value = 43.409393099401974
$ws.Range("AS3").Value = 47.289831801062746
$ws.Range("AT3").Value = 60.437187185095432
$ws.Range("AU3").Value = 25.824627668321781
$ws.Range("AV3").Value = 48.175406816009556
$ws.Range("AW3").Value = 48.471267249543352
$ws.Range("AX3").Value = 37.756105655407936
$ws.Range("AY3").Value = 67.278804168408641
